$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.946.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.42%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.893.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.03%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7748'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.47%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.41%  '

# Row 7
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3137'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.92%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.68'
$ws.Range('D9').Style = 'Normal'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07351'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.27%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08061'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.19%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7734'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.68%  '

# Row 13
$ws.Range('E13').Value = '  +2.91%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.862.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.05%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.13'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.78%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.229'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.65%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.890.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.20%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.02'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.09%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.85%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007867'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.152'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.59%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9998'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.02%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.118.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.28%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '

# Row 25
$ws.Range('E25').Value = '  -4.16%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.446'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.99%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.20%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.38%  '

# Row 29
$ws.Range('E29').Value = '  -1.42%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.424'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.73%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.541'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.478'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.84%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05568'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.08%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.067'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.68%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.237'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.81%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7497'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.60%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9995'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.10%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.684'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.13%  '

# Row 39
$ws.Range('E39').Value = '  +1.07%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.789'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.09%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.71%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.100.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.43%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.006'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.34%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8504'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.09%  '

# Row 46
$ws.Range('E46').Value = '  +0.01%  '

# Row 47
$ws.Range('E47').Value = '  +1.00%  '

# Row 48
$ws.Range('E48').Value = '  -0.62%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.547'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.68%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.765'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.08%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.991'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.17%  '
